$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shorten the "Address" (column C) text for a set of rows: drop the
# building/street-level portion, keeping just the city/state/country tail.
$ws.Range("C2").Value  = " Kannapiran Mill Road Sowripalayam, Coimbatore, Tamilnadu India"
$ws.Range("C3").Value  = " Mahabalipuram, New Delhi - 110074"
$ws.Range("C4").Value  = " POLLACHI, COIMBATORE – 642004."
$ws.Range("C5").Value  = "Vadodara, Gujarat 389390, India"
$ws.Range("C6").Value  = " New Delhi, Delhi 110068"
$ws.Range("C7").Value  = " Thondamuthur, Tamil Nadu 641109"
$ws.Range("C10").Value = " Amritsar, Punjab 143001"
$ws.Range("C12").Value = "Palakkad, Kerala,India"
$ws.Range("C13").Value = " Gurgaon - 122001, India"
$ws.Range("C14").Value = "Mulshi, Maharashtra 412115, India"
$ws.Range("C15").Value = "Faridabad, Haryana 121101, India"

$ws.Range("C15").Select()

$wb.Save()
